# Apply the edits described by the commit:
# "removed ER tags from non-ER templates and non-ER tags"
#
# 1. Rename the "SwateTemplateMetadata" sheet to "isa_template"
# 2. Clear the ER / ER Term Accession Number / ER Term Source REF values
#    (these referenced the "PRIDE" / DPBO ontology entries which do not
#    belong to this non-ER template)
# 3. Clear the 4th Tags / Tags Term Accession Number / Tags Term Source REF
#    column (the "PRIDE" tag + DPBO term that should not be tagged here)
# 4. Update the active selection on that sheet

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Name = "isa_template"

# Remove the ER row values (row 8-10, column B)
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()

# Remove the 4th Tags column values (column E, rows 12-14)
$ws.Range("E12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("E14").ClearContents()

# The wrapped-text hyperlink cell in row 13 no longer needs extra height
# once the sheet is resaved; let Excel recompute the row height.
$ws.Rows.Item(13).AutoFit() | Out-Null

# Update selection to match the recorded view state
$ws.Activate()
$ws.Range("B18").Select() | Out-Null
